$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "744"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1695687.79"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1021"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3631453.47"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "662"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2143408.78"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "365"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1364308.18"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "170"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "614976.34"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "222"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "600362.00"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "289"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "747727.45"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "562"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2326182.70"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "390"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1332198.57"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "16"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63500.00"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "166"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "446117.22"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "87"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "417409.98"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "131"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "535972.25"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "4"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13109.00"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "386"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1041140.74"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "625"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2482013.99"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "426"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1486861.40"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "3816"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8800141.47"
$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "4055"
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value = "14275525.38"
$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "4189"
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = "13308255.27"
$ws.Range("C58").NumberFormat = "@"
$ws.Range("C58").Value = "90"
$ws.Range("D58").NumberFormat = "@"
$ws.Range("D58").Value = "328110.47"
$ws.Range("C76").NumberFormat = "@"
$ws.Range("C76").Value = "929"
$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value = "3247240.26"
$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "524"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "1726125.47"
